# Update "想去人数" (F column) counts on both the "展览" and "全部类型" sheets.
# Each entry is identified by its row number on the respective sheet.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet4 = $wb.Worksheets.Item("全部类型")

# Row => new F value, for sheet "展览"
$updates1 = @{
    2  = 2717
    4  = 356
    6  = 1141
    13 = 9210
    20 = 633
    25 = 2191
    27 = 1892
    31 = 280
    37 = 297
    41 = 683
    42 = 42
    43 = 1399
    47 = 647
}

foreach ($row in $updates1.Keys) {
    $sheet1.Range("F$row").Value = $updates1[$row]
}

# Row => new F value, for sheet "全部类型"
$updates4 = @{
    2  = 2717
    3  = 356
    6  = 1141
    10 = 9210
    19 = 633
    22 = 2191
    23 = 1892
    26 = 280
    32 = 297
    39 = 684
    41 = 42
    42 = 1399
    47 = 647
}

foreach ($row in $updates4.Keys) {
    $sheet4.Range("F$row").Value = $updates4[$row]
}
